$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add materials for session 05 (row 6): Folien (slides) and Aufgaben (exercises)
$ws.Range("E6").Value = "slides/slides.html#/sitzung-05-open-science"
$ws.Range("F6").Value = "exercises/e05.html"

# Update the active cell selection to reflect where the author ended up
$ws.Range("F7").Select()
